# This script applies the edit described in the commit:
# "feat: add 2022-Q3 data"
#
# It performs two logical changes to the workbook:
#   1. Inserts a new summary row for 2022-Q3 at the top of the data in the
#      "Total" (sheet 1) worksheet, shifting the existing quarters down.
#   2. Inserts a brand-new worksheet named "2022-Q3" (as the 2nd tab, right
#      after the totals sheet) containing the per-fund holdings detail for
#      that quarter, built from the existing Q2 sheet as a style template.

$wb = $excel.ActiveWorkbook

# ======================================================================
# Step 1: Update the totals worksheet (tab 1) with the new Q3 summary row
# ======================================================================
$wsTotal = $wb.Worksheets.Item(1)

$totalsData = @(
    @(0, '2022-Q3', 43, 9.95),
    @(1, '2022-Q2', 27, 5.66),
    @(2, '2022-Q1', 15, 4.81),
    @(3, '2021-Q4', 16, 4.37),
    @(4, '2021-Q3', 10, 2.84),
    @(5, '2021-Q2', 3, 1.58)
)

# The bold + thin-border index style used on column A (A2:A6) needs to be
# extended onto the new row that is being appended (old sheet only went to A6).
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A7").PasteSpecial(-4122)
$wsTotal.Application.CutCopyMode = $false

for ($i = 0; $i -lt $totalsData.Length; $i++) {
    $r = $i + 2
    $d = $totalsData[$i]
    $wsTotal.Cells.Item($r, 1).Value = $d[0]
    $wsTotal.Cells.Item($r, 2).Value = $d[1]
    $wsTotal.Cells.Item($r, 3).Value = $d[2]
    $wsTotal.Cells.Item($r, 4).Value = $d[3]
}

# ======================================================================
# Step 2: Create the new "2022-Q3" worksheet as a copy of "2022-Q2" (tab 2)
# so that it inherits identical column layout / header styling, then
# overwrite its contents with the Q3 fund holdings.
# ======================================================================
$wsQ2Template = $wb.Worksheets.Item(2)

# Copy-before-itself: the duplicate is inserted directly in front of the
# original "2022-Q2" sheet, i.e. it becomes the new tab 2, and every sheet
# that used to follow (Q2, Q1, 2021-Q4, 2021-Q3, 2021-Q2) shifts one slot right.
$wsQ2Template.Copy($wsQ2Template)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# The template sheet only has 28 rows (header + 27 funds); Q3 needs 44 rows
# (header + 43 funds), so first extend the header-style formatting of the
# template's last templated row down across the additional rows that will
# be needed, copying the whole row (A:H) so column A keeps its bold/border
# style and text columns keep a plain format before we force them to Text.
$templateLastDataRow = 28
$newLastDataRow = 44
$wsQ3.Range("A2:H2").Copy()
$wsQ3.Range("A" + ($templateLastDataRow + 1) + ":H" + $newLastDataRow).PasteSpecial(-4122)
$wsQ3.Application.CutCopyMode = $false

# Fund code (B), fund name (C), fund size / positions / weight / value (D:G)
# are all stored as plain text in this workbook (fund codes have leading
# zeroes, and D:G keep a fixed number of decimals) -- format the full data
# range as Text up front so assigning numeric-looking strings does not get
# auto-converted into numbers by Excel.
$wsQ3.Range("B2:G" + $newLastDataRow).NumberFormat = "@"

$q3Data = @(
    @(0, '005176', '富国精准医疗灵活配置混合', '33.00', '90.96', '4.30', '1.4190', 9),
    @(1, '009664', '汇添富医疗积极成长一年持有期混合A', '29.13', '66.68', '4.34', '1.2642', 6),
    @(2, '007553', '中信建投医改灵活配置混合C', '11.66', '94.99', '9.10', '1.0611', 6),
    @(3, '002408', '中信建投医改灵活配置混合A', '10.95', '94.99', '9.10', '0.9964', 6),
    @(4, '002708', '大摩健康产业混合A', '21.51', '92.37', '4.18', '0.8991', 9),
    @(5, '000711', '嘉实医疗保健股票', '14.17', '94.06', '5.87', '0.8318', 4),
    @(6, '005303', '嘉实医药健康股票A', '11.55', '94.88', '5.15', '0.5948', 3),
    @(7, '011868', '中信建投远见回报混合A', '6.14', '94.99', '5.79', '0.3555', 6),
    @(8, '010090', '中信建投医药健康混合A', '2.88', '94.86', '9.30', '0.2678', 6),
    @(9, '001915', '宝盈医疗健康沪港深股票', '5.15', '91.59', '4.97', '0.2560', 6),
    @(10, '070022', '嘉实领先成长混合', '4.99', '86.27', '4.73', '0.2360', 9),
    @(11, '005304', '嘉实医药健康股票C', '4.19', '94.88', '5.15', '0.2158', 3),
    @(12, '010091', '中信建投医药健康混合C', '1.98', '94.86', '9.30', '0.1841', 6),
    @(13, '014030', '大摩健康产业混合C', '4.27', '92.37', '4.18', '0.1785', 9),
    @(14, '009665', '汇添富医疗积极成长一年持有期混合C', '3.63', '66.68', '4.34', '0.1575', 6),
    @(15, '002300', '长盛医疗行业量化配置股票', '2.39', '93.48', '4.91', '0.1173', 9),
    @(16, '014867', '摩根士丹利华鑫优悦安和混合C', '1.32', '93.41', '8.04', '0.1061', 5),
    @(17, '000523', '国投瑞银医疗保健混合A', '1.95', '92.21', '4.75', '0.0926', 8),
    @(18, '013357', '大摩沪港深精选混合C', '1.53', '92.27', '5.52', '0.0845', 8),
    @(19, '011869', '中信建投远见回报混合C', '1.39', '94.99', '5.79', '0.0805', 6),
    @(20, '009893', '摩根士丹利华鑫优悦安和混合A', '0.87', '93.41', '8.04', '0.0699', 5),
    @(21, '000684', '长盛养老健康产业灵活配置混合', '1.35', '92.15', '4.60', '0.0621', 9),
    @(22, '970023', '天风天盈一年定期开放混合', '1.24', '70.03', '4.28', '0.0531', 9),
    @(23, '001056', '北信瑞丰健康生活主题灵活配置混合', '0.97', '93.32', '4.31', '0.0418', 8),
    @(24, '005108', '圆信永丰双利优选定期开放灵活配置混合', '0.61', '91.10', '6.56', '0.0400', 1),
    @(25, '013356', '大摩沪港深精选混合A', '0.68', '92.27', '5.52', '0.0375', 8),
    @(26, '014220', '恒越医疗健康精选混合A', '0.72', '88.76', '5.12', '0.0369', 3),
    @(27, '001965', '圆信永丰兴源灵活配置混合A', '0.48', '93.86', '6.99', '0.0336', 2),
    @(28, '008412', '长盛竞争优势股票A', '0.72', '91.35', '4.53', '0.0326', 9),
    @(29, '015032', '中融医药消费混合A', '0.54', '90.81', '5.05', '0.0273', 3),
    @(30, '005520', '国投瑞银创新医疗混合', '0.41', '93.90', '4.49', '0.0184', 8),
    @(31, '001966', '圆信永丰兴源灵活配置混合C', '0.26', '93.86', '6.99', '0.0182', 2),
    @(32, '008413', '长盛竞争优势股票C', '0.39', '91.35', '4.53', '0.0177', 9),
    @(33, '006241', '中融医疗健康精选混合C', '0.46', '88.49', '3.24', '0.0149', 9),
    @(34, '014221', '恒越医疗健康精选混合C', '0.29', '88.76', '5.12', '0.0148', 3),
    @(35, '006274', '圆信永丰医药健康混合', '0.16', '93.60', '6.93', '0.0111', 1),
    @(36, '350007', '天治趋势精选混合', '0.37', '82.33', '2.05', '0.0076', 8),
    @(37, '001563', '华富健康文娱灵活配置混合', '0.13', '93.80', '4.80', '0.0062', 1),
    @(38, '011082', '国投瑞银医疗保健混合C', '0.11', '92.21', '4.75', '0.0052', 8),
    @(39, '006240', '中融医疗健康精选混合A', '0.05', '88.49', '3.24', '0.0016', 9),
    @(40, '004724', '先锋聚元灵活配置混合A', '0.04', '94.36', '2.47', '0.0010', 10),
    @(41, '004725', '先锋聚元灵活配置混合C', '0.04', '94.36', '2.47', '0.0010', 10),
    @(42, '015033', '中融医药消费混合C', '0.02', '90.81', '5.05', '0.0010', 3)
)

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $r = $i + 2
    $d = $q3Data[$i]
    $wsQ3.Cells.Item($r, 1).Value = $d[0]
    $wsQ3.Cells.Item($r, 2).Value = $d[1]
    $wsQ3.Cells.Item($r, 3).Value = $d[2]
    $wsQ3.Cells.Item($r, 4).Value = $d[3]
    $wsQ3.Cells.Item($r, 5).Value = $d[4]
    $wsQ3.Cells.Item($r, 6).Value = $d[5]
    $wsQ3.Cells.Item($r, 7).Value = $d[6]
    $wsQ3.Cells.Item($r, 8).Value = $d[7]
}

Write-Output "2022-Q3 sheet and totals updated successfully"
